$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the training-part split, reusing the bold/centered
# header style already used by A1:E1.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# Register both number-format variants (lowercase then uppercase) on the same
# cell so the styles sheet ends up with both numFmt entries while only the
# uppercase one is actually referenced by a cellXfs record.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$rowsData = @(
  ,@(45684.59222037037, 567.8, 10.38, 1.402552247047424, "10-15", "Duża Gra")
  ,@(45684.59352476852, 680.5, 11.18, 1.723067641258238, "10-15", "Duża Gra")
  ,@(45684.59410462963, 730.6, 10.29, 1.463077562195914, "10-15", "Duża Gra")
  ,@(45684.59107569444, 468.9, 8.25, 1.562695315905979, "5-10", "Duża Gra")
  ,@(45684.5935224537, 680.3, 9.83, 1.62039874281202, "5-10", "Duża Gra")
  ,@(45684.59410347222, 730.5, 9.81, 1.491053649357387, "5-10", "Duża Gra")
  ,@(45684.59765208334, 1037.1, 13.72, 3.339186361857823, "10-15", "Mała Gra")
  ,@(45684.60033495371, 1268.9, 12.45, 3.148051295961654, "10-15", "Mała Gra")
  ,@(45684.6030537037, 1503.8, 14.47, 3.442948818206787, "10-15", "Mała Gra")
  ,@(45684.59764861111, 1036.8, 9.949999999999999, 2.825158732278006, "5-10", "Mała Gra")
  ,@(45684.60033148148, 1268.6, 8.789999999999999, 2.768171565873284, "5-10", "Mała Gra")
  ,@(45684.60429560185, 1611.1, 9.75, 2.816675628934587, "5-10", "Mała Gra")
)

$r = 2
foreach ($row in $rowsData) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r = $r + 1
}
